$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '19.975.41'
$ws.Range('E2').Value = '  -7.89%  '
$ws.Range('D3').Value = '1.414.70'
$ws.Range('E3').Value = '  -7.95%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9988'
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9997'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '274.96'
$ws.Range('E6').Value = '  -5.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3689'
$ws.Range('E7').Value = '  -6.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3128'
$ws.Range('E8').Value = '  -2.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '39.61'
$ws.Range('E9').Value = '  -8.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.039'
$ws.Range('E10').Value = '  -3.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06499'
$ws.Range('E11').Value = '  -9.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9989'
$ws.Range('E12').Value = '  -0.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.467'
$ws.Range('E13').Value = '  -5.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '17.68'
$ws.Range('E14').Value = '  -3.98%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.199'
$ws.Range('E15').Value = '  -6.40%  '
$ws.Range('D16').Value = '1.411.47'
$ws.Range('E16').Value = '  -8.43%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001017'
$ws.Range('E17').Value = '  -7.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.05684'
$ws.Range('E18').Value = '  -14.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.000'
$ws.Range('E19').Value = '  +0.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.68'
$ws.Range('E20').Value = '  -15.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.607'
$ws.Range('E21').Value = '  -8.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.76'
$ws.Range('E22').Value = '  -4.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.07'
$ws.Range('E23').Value = '  +2.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.290'
$ws.Range('E24').Value = '  -3.56%  '
$ws.Range('D25').Value = '19.986.80'
$ws.Range('E25').Value = '  -7.82%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.270'
$ws.Range('E26').Value = '  -4.92%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '135.40'
$ws.Range('E27').Value = '  -10.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.01'
$ws.Range('E28').Value = '  -7.97%  '
$ws.Range('D29').Value = '1.568.02'
$ws.Range('E29').Value = '  -8.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '109.68'
$ws.Range('E30').Value = '  -6.70%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.094'
$ws.Range('E31').Value = '  -16.43%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.357'
$ws.Range('E32').Value = '  -11.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.8269'
$ws.Range('E33').Value = '  -14.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07693'
$ws.Range('E34').Value = '  -4.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '8.404'
$ws.Range('E35').Value = '  -1.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.455'
$ws.Range('E36').Value = '  -2.87%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.893'
$ws.Range('E37').Value = '  -5.73%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05836'
$ws.Range('E38').Value = '  -2.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9994'
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02071'
$ws.Range('E40').Value = '  -6.77%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '10.51'
$ws.Range('E41').Value = '  -6.99%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1916'
$ws.Range('E42').Value = '  -6.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5312'
$ws.Range('E44').Value = '  -8.64%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.31'
$ws.Range('E45').Value = '  -6.81%  '
$ws.Range('E46').Value = '  -5.32%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5174'
$ws.Range('E47').Value = '  -7.30%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '112.76'
$ws.Range('E48').Value = '  -2.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.775'
$ws.Range('E49').Value = '  -6.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.043'
$ws.Range('E50').Value = '  -9.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9986'
$ws.Range('E51').Value = '  -0.25%  '
